$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Textfeld 1" shape that lists the authors, addressed by name
# so this keeps working even if shape ordering ever changes.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Textfeld 1") {
        $shp = $candidate
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(2)
}

$tr = $shp.TextFrame.TextRange

# Find the paragraph that currently reads "Denise Langhals" (two runs:
# "Denise " + "Langhals") and replace it with a single run reading
# "Denise Langhof", merging the previously split runs.
$paraCount = $tr.Paragraphs().Count
$target = $null
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "Denise*") {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    $target = $tr.Paragraphs(4, 1)
}

$full = $target.Characters(1, $target.Length)
$full.Text = "Denise Langhof"
